$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text formatting
# (values such as "1.00" or "51.286.85" must remain literal text, not numbers)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "51.286.85"
$ws.Range("D3").Value = "2.976.49"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "380.86"
$ws.Range("E5").Value = "  +1.84%  "
$ws.Range("D6").Value = "102.89"
$ws.Range("E6").Value = "  +2.25%  "
$ws.Range("D7").Value = "0.545"
$ws.Range("E7").Value = "  +2.14%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.591"
$ws.Range("E9").Value = "  +1.75%  "
$ws.Range("D10").Value = "36.62"
$ws.Range("E10").Value = "  +1.55%  "
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("D12").Value = "0.0859"
$ws.Range("E12").Value = "  +1.72%  "
$ws.Range("D13").Value = "3.452.26"
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("D14").Value = "7.79"
$ws.Range("E14").Value = "  +4.16%  "
$ws.Range("D15").Value = "18.38"
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("D16").Value = "2.998.01"
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("D17").Value = "11.20"
$ws.Range("E17").Value = "  +1.61%  "
$ws.Range("D18").Value = "0.998"
$ws.Range("E18").Value = "  +2.94%  "
$ws.Range("D19").Value = "51.361.32"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("D21").Value = "12.54"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").Value = "70.26"
$ws.Range("E23").Value = "  +2.45%  "
$ws.Range("D24").Value = "267.14"
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("D25").Value = "3.21"
$ws.Range("E25").Value = "  +3.27%  "
$ws.Range("D26").Value = "7.81"
$ws.Range("E26").Value = "  -3.07%  "
$ws.Range("E27").Value = "  +1.99%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").Value = "26.02"
$ws.Range("E29").Value = "  +1.79%  "
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("E31").Value = "  -1.48%  "
$ws.Range("D32").Value = "10.32"
$ws.Range("E32").Value = "  +4.12%  "
$ws.Range("E33").Value = "  +5.03%  "
$ws.Range("D34").Value = "51.48"
$ws.Range("E35").Value = "  +1.48%  "
$ws.Range("E36").Value = "  -0.70%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "3.25"
$ws.Range("E38").Value = "  +3.95%  "
$ws.Range("E39").Value = "  +1.67%  "
$ws.Range("D40").Value = "16.71"
$ws.Range("E40").Value = "  +2.62%  "
$ws.Range("E41").Value = "  +3.17%  "
$ws.Range("E42").Value = "  +3.12%  "
$ws.Range("D43").Value = "124.47"
$ws.Range("E43").Value = "  +4.00%  "
$ws.Range("E44").Value = "  +10.91%  "
$ws.Range("D45").Value = "21.64"
$ws.Range("E45").Value = "  +2.86%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("E47").Value = "  +4.30%  "
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("D49").Value = "2.034.13"
$ws.Range("E49").Value = "  +2.40%  "
$ws.Range("E50").Value = "  +1.52%  "
$ws.Range("D51").Value = "0.533"
$ws.Range("E51").Value = "  +16.26%  "
